$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holds price-like text values (e.g. "245.80"). Excel would normally
# auto-coerce such text into a Number when assigned via .Value, so we temporarily
# format the target price cells as Text, write the values, then restore the
# default "Normal" cell style so no stray number-format / style is left behind.
$priceCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D22","D23","D26","D27","D28","D40","D41","D42","D43","D44","D45","D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "245.83"
$ws.Range("D3").Value = "25.40"
$ws.Range("D4").Value = "5.144"
$ws.Range("D5").Value = "0.05596"
$ws.Range("D7").Value = "3.012"
$ws.Range("D8").Value = "0.8174"
$ws.Range("D9").Value = "0.8420"
$ws.Range("D10").Value = "0.1341"
$ws.Range("D11").Value = "0.03172"
$ws.Range("D12").Value = "0.02849"
$ws.Range("D13").Value = "0.09395"
$ws.Range("D14").Value = "0.001522"
$ws.Range("D15").Value = "0.0005941"
$ws.Range("E15").Value = "14OneONEWorstin24h"
$ws.Range("D16").Value = "0.006246"
$ws.Range("D17").Value = "3.512"
$ws.Range("D18").Value = "2.082"
$ws.Range("D20").Value = "0.06955"
$ws.Range("D22").Value = "3.747"
$ws.Range("D23").Value = "0.04730"
$ws.Range("D26").Value = "0.004274"
$ws.Range("D27").Value = "0.00009699"
$ws.Range("D28").Value = "0.0001388"
$ws.Range("D40").Value = "0.03656"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006220"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1054"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002636"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "0.008383"
$ws.Range("D45").Value = "0.00005304"
$ws.Range("D48").Value = "0.002284"

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
